$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row to lowercase column names
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "age"
$ws.Range("C1").Value = "height"

# Select cell C1 as the active cell (matches the sheetView selection in the diff)
$ws.Range("C1").Select()
